$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.118.64'
$ws.Range('D2').Style = $origStyle
$origStyle = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('E2').Style = $origStyle
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.900.50'
$ws.Range('D3').Style = $origStyle
$origStyle = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E3').Style = $origStyle
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.97'
$ws.Range('D5').Style = $origStyle
$origStyle = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.19%  '
$ws.Range('E5').Style = $origStyle
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.690'
$ws.Range('D6').Style = $origStyle
$origStyle = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('E6').Style = $origStyle
$origStyle = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('E7').Style = $origStyle
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.28'
$ws.Range('D8').Style = $origStyle
$origStyle = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.48%  '
$ws.Range('E8').Style = $origStyle
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.357'
$ws.Range('D9').Style = $origStyle
$origStyle = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.80%  '
$ws.Range('E9').Style = $origStyle
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.71'
$ws.Range('D10').Style = $origStyle
$origStyle = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E10').Style = $origStyle
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0748'
$ws.Range('D11').Style = $origStyle
$origStyle = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.71%  '
$ws.Range('E11').Style = $origStyle
$origStyle = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.20%  '
$ws.Range('E12').Style = $origStyle
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.177.75'
$ws.Range('D13').Style = $origStyle
$origStyle = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('E13').Style = $origStyle
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.98'
$ws.Range('D14').Style = $origStyle
$origStyle = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.42%  '
$ws.Range('E14').Style = $origStyle
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.732'
$ws.Range('D15').Style = $origStyle
$origStyle = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.76%  '
$ws.Range('E15').Style = $origStyle
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.96'
$ws.Range('D16').Style = $origStyle
$origStyle = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.09%  '
$ws.Range('E16').Style = $origStyle
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.912.65'
$ws.Range('D17').Style = $origStyle
$origStyle = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('E17').Style = $origStyle
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.138.70'
$ws.Range('D18').Style = $origStyle
$origStyle = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('E18').Style = $origStyle
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.40'
$ws.Range('D19').Style = $origStyle
$origStyle = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('E19').Style = $origStyle
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0829'
$ws.Range('D20').Style = $origStyle
$origStyle = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('E20').Style = $origStyle
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.31'
$ws.Range('D21').Style = $origStyle
$origStyle = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('E21').Style = $origStyle
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.94'
$ws.Range('D22').Style = $origStyle
$origStyle = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.61%  '
$ws.Range('E22').Style = $origStyle
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.03'
$ws.Range('D23').Style = $origStyle
$origStyle = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.97%  '
$ws.Range('E23').Style = $origStyle
$origStyle = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('E24').Style = $origStyle
$origStyle = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.88%  '
$ws.Range('E25').Style = $origStyle
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = $origStyle
$origStyle = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.81%  '
$ws.Range('E26').Style = $origStyle
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.67'
$ws.Range('D27').Style = $origStyle
$origStyle = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('E27').Style = $origStyle
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.53'
$ws.Range('D28').Style = $origStyle
$origStyle = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E28').Style = $origStyle
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.48'
$ws.Range('D29').Style = $origStyle
$origStyle = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('E29').Style = $origStyle
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.129'
$ws.Range('D30').Style = $origStyle
$origStyle = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('E30').Style = $origStyle
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.128.07'
$ws.Range('D31').Style = $origStyle
$origStyle = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E31').Style = $origStyle
$origStyle = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +12.62%  '
$ws.Range('E32').Style = $origStyle
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0606'
$ws.Range('D33').Style = $origStyle
$origStyle = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.77%  '
$ws.Range('E33').Style = $origStyle
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.30'
$ws.Range('D34').Style = $origStyle
$origStyle = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.83%  '
$ws.Range('E34').Style = $origStyle
$origStyle = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.46%  '
$ws.Range('E35').Style = $origStyle
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.19'
$ws.Range('D36').Style = $origStyle
$origStyle = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('E36').Style = $origStyle
$origStyle = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E37').Style = $origStyle
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.853'
$ws.Range('D38').Style = $origStyle
$origStyle = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.74%  '
$ws.Range('E38').Style = $origStyle
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.99'
$ws.Range('D39').Style = $origStyle
$origStyle = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('E39').Style = $origStyle
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '102.50'
$ws.Range('D40').Style = $origStyle
$origStyle = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +14.47%  '
$ws.Range('E40').Style = $origStyle
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.27'
$ws.Range('D41').Style = $origStyle
$origStyle = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.06%  '
$ws.Range('E41').Style = $origStyle
$origStyle = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('E42').Style = $origStyle
$origStyle = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('E44').Style = $origStyle
$origStyle = $ws.Range('B45').Style
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Maker'
$ws.Range('B45').Style = $origStyle
$origStyle = $ws.Range('C45').Style
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C45').Style = $origStyle
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.317.94'
$ws.Range('D45').Style = $origStyle
$origStyle = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E45').Style = $origStyle
$origStyle = $ws.Range('B46').Style
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('B46').Style = $origStyle
$origStyle = $ws.Range('C46').Style
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C46').Style = $origStyle
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.41'
$ws.Range('D46').Style = $origStyle
$origStyle = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('E46').Style = $origStyle
$origStyle = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E47').Style = $origStyle
$origStyle = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('E48').Style = $origStyle
$origStyle = $ws.Range('B49').Style
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Gas'
$ws.Range('B49').Style = $origStyle
$origStyle = $ws.Range('C49').Style
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('C49').Style = $origStyle
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.21'
$ws.Range('D49').Style = $origStyle
$origStyle = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('E49').Style = $origStyle
$origStyle = $ws.Range('B50').Style
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('B50').Style = $origStyle
$origStyle = $ws.Range('C50').Style
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C50').Style = $origStyle
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.58'
$ws.Range('D50').Style = $origStyle
$origStyle = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('E50').Style = $origStyle
$origStyle = $ws.Range('B51').Style
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('B51').Style = $origStyle
$origStyle = $ws.Range('C51').Style
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C51').Style = $origStyle
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0745'
$ws.Range('D51').Style = $origStyle
$origStyle = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.79%  '
$ws.Range('E51').Style = $origStyle
